$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.975.79"
$ws.Range("E2").Value = "  +1.54%  "
$ws.Range("D3").Value = "3.504.81"
$ws.Range("E3").Value = "  +0.24%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "601.85"
$ws.Range("E5").Value = "  +2.22%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.15"
$ws.Range("E6").Value = "  +1.85%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.608"
$ws.Range("E7").Value = "  +0.36%  "
$ws.Range("D8").Value = "3.497.66"
$ws.Range("E8").Value = "  +0.23%  "
$ws.Range("E9").Value = "  +0.11%  "
$ws.Range("E10").Value = "  -0.85%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.25"
$ws.Range("E11").Value = "  +7.31%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.580"
$ws.Range("E12").Value = "  +1.03%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "46.22"
$ws.Range("E13").Value = "  -1.31%  "
$ws.Range("E14").Value = "  -0.58%  "
$ws.Range("D15").Value = "4.066.63"
$ws.Range("E15").Value = "  +0.45%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "8.27"
$ws.Range("E16").Value = "  -0.98%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "604.67"
$ws.Range("E17").Value = "  -2.25%  "
$ws.Range("D18").Value = "70.040.67"
$ws.Range("E18").Value = "  +1.76%  "
$ws.Range("D19").Value = "3.491.94"
$ws.Range("E19").Value = "  +0.30%  "
$ws.Range("E20").Value = "  +0.86%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.11"
$ws.Range("E21").Value = "  -1.03%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.867"
$ws.Range("E22").Value = "  -1.00%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.10"
$ws.Range("E23").Value = "  -18.44%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "15.44"
$ws.Range("E24").Value = "  -2.47%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "95.39"
$ws.Range("E25").Value = "  -0.57%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.70"
$ws.Range("E26").Value = "  -2.06%  "
$ws.Range("E27").Value = "  -0.05%  "
$ws.Range("E28").Value = "  -1.55%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.98"
$ws.Range("E29").Value = "  +2.10%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.95"
$ws.Range("E30").Value = "  -2.56%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "707.82"
$ws.Range("E31").Value = "  +22.61%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.01"
$ws.Range("E32").Value = "  -3.03%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "8.07"
$ws.Range("E33").Value = "  -4.45%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.90"
$ws.Range("E34").Value = "  +0.57%  "
$ws.Range("E35").Value = "  -2.53%  "
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0996"
$ws.Range("E36").Value = "  -1.63%  "
$ws.Range("B37").Value = "dogwifhat"
$ws.Range("C37").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.54"
$ws.Range("E37").Value = "  +0.29%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "10.66"
$ws.Range("E38").Value = "  -0.40%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0470"
$ws.Range("E39").Value = "  +7.74%  "
$ws.Range("B40").Value = "FirstDigitalUSD"
$ws.Range("C40").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("E40").Value = "  +0.23%  "
$ws.Range("B41").Value = "OKB"
$ws.Range("C41").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "56.42"
$ws.Range("E41").Value = "  -1.17%  "
$ws.Range("E42").Value = "  +4.00%  "
$ws.Range("D43").Value = "3.324.18"
$ws.Range("E43").Value = "  -2.64%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.311"
$ws.Range("E44").Value = "  -3.98%  "
$ws.Range("E45").Value = "  +4.27%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "32.16"
$ws.Range("E46").Value = "  -1.75%  "
$ws.Range("D47").Value = "0.0₃0687"
$ws.Range("E47").Value = "  -1.29%  "
$ws.Range("E48").Value = "  -0.63%  "
$ws.Range("E49").Value = "  +0.76%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "133.24"
$ws.Range("E50").Value = "  +0.98%  "

# Restore default style on cells where we temporarily forced text format,
# so no stray number-format style lingers on the saved workbook.
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D50").Style = "Normal"
